# Apply "added my project tracker" updates to the Tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracker")

# --- Row 5 / Row 6: split the old combined "serial comms and memory
# management scheme" whiteboard goal into two separate goals ---
$ws.Range("C5").Value = "Whiteboard/brainstorm serial comms implementation."

$ws.Range("C6").Value = "Whiteboard/brainstorm  memory management scheme."
# C6 had no fill/border formatting yet - match the styling already used by
# the other goal cells in column C (C3:C5) by copying their format over.
$ws.Range("C5").Copy() | Out-Null
$ws.Range("C6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Row 9: update the weekly reflection text for week three (column C) ---
$ws.Range("C9").Value = "All goals except PID code were completed. Noah and I were still finishing the I2C lab so most of the PIC coding was done by Jacob. I led a coule of brainstorming sessions to ensure that the team is aligned on how our firmware will control the hardware of the system. For PID control I got C code for PID math and now I need to implement it in PIC16F1788."
# The extra sentence makes the wrapped text taller - match Excel's autofit result.
$ws.Rows.Item(9).RowHeight = 115.2

# --- New "Resources" links under the tracker (rows 13-15) ---
$ws.Range("A13").Value = "https://github.com/Christopher-isu/Farm-Bureau-Project/blob/main/IO_Block_diagram_v1.svg"
$ws.Range("C13").Value = "https://github.com/Christopher-isu/Farm-Bureau-Project/blob/main/Docs/Whiteboard%20Sessions/IOscheme.jpg"
$ws.Range("C14").Value = "https://github.com/Christopher-isu/Farm-Bureau-Project/blob/main/Docs/Whiteboard%20Sessions/MCU-UI-commStructure.jpg"
$ws.Range("C15").Value = "https://github.com/Christopher-isu/Farm-Bureau-Project/blob/main/Docs/Whiteboard%20Sessions/process_cycle.jpg"

# Give the new rows the same row height / wrap-text look as the rest of the sheet.
$ws.Rows.Item(13).RowHeight = 43.2
$ws.Rows.Item(14).RowHeight = 43.2
$ws.Rows.Item(15).RowHeight = 43.2

$ws.Range("C13:C15").WrapText = $true
$ws.Range("C13:C15").Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$ws.Range("C13:C15").Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous

# Turn A13 into a real hyperlink pointing at the same URL shown as its text.
$ws.Hyperlinks.Add($ws.Range("A13"), "https://github.com/Christopher-isu/Farm-Bureau-Project/blob/main/IO_Block_diagram_v1.svg")

# Trailing blank rows below the new Resources block, matching the sheet's
# natural (non-custom) row height.
$ws.Rows.Item(16).RowHeight = 14.4
$ws.Rows.Item(17).RowHeight = 14.4
$ws.Rows.Item(18).RowHeight = 14.4
$ws.Rows.Item(19).RowHeight = 14.4
$ws.Rows.Item(20).RowHeight = 14.4
$ws.Rows.Item(21).RowHeight = 14.4

# Match the final selection left behind in the saved workbook.
$ws.Range("D9").Select()
